# Daily attendance processing - 2026-01-24 11:33:47
# Swap the order of "System" and the email address in the "Recorded By"
# column (G) wherever the value is "System, dnasr281@gmail.com", turning
# it into "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    if ($cell.Value() -eq $oldValue) {
        $cell.Value = $newValue
    }
}
